# Update table1.xlsx (Logreg_day1.ipynb output) to reflect the new set of
# engineered features: the old "hadMeasurmentDayOne" row is replaced by two
# derived rows ("hadMeasurmentDayOne_chart" and "hadMeasurmentDayOne_lab")
# that are appended at the end of the table, together with a brand new
# "measurment_before" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 14 ("hadMeasurmentDayOne, n (%)") needs to move to the bottom
# of the table (rows 20 & 21, duplicated), while rows 15-20 shift up to
# become rows 14-19. Duplicate the row (values + formatting) one row below
# its eventual target first, so that deleting row 14 afterwards shifts
# everything (including our duplicates) up by exactly one row, landing them
# on rows 20 and 21.
$ws.Range("A14:H14").Copy($ws.Range("A21:H21"))
$ws.Range("A14:H14").Copy($ws.Range("A22:H22"))

# Remove the original row; rows 15-22 shift up to 14-21.
$ws.Rows(14).Delete()

# Rename the two duplicated rows (now at 20 and 21).
$ws.Range("A20").Value = "hadMeasurmentDayOne_chart, n (%)"
$ws.Range("A21").Value = "hadMeasurmentDayOne_lab, n (%)"

# Build the brand-new row 22 ("measurment_before, n (%)"). Reuse the
# formatting (and, for column B, the literal "True" text) from row 21 so
# that the values keep the same text typing/style as the rest of the table
# instead of being auto-coerced into booleans/numbers.
$ws.Range("B21").Copy($ws.Range("B22"))
$ws.Range("A21").Copy($ws.Range("A22"))
$ws.Range("A22").Value = "measurment_before, n (%)"

$ws.Range("C22").Value = 23326
$ws.Range("D22").Value = "1601 (100.0)"
$ws.Range("E22").Value = "68 (100.0)"
$ws.Range("F22").Value = "253 (100.0)"
$ws.Range("G22").Value = "112 (100.0)"
$ws.Range("H22").Value = "1168 (100.0)"

Write-Output "table1.xlsx rows updated"
